# "Updated grades for HW 7" -- Grades.xlsx
#
# Adds a new "HW 8 / 7 CPP" grade group (3 new blank group slots total,
# only the first one gets data) and fills in previously-blank HW5/HW6/HW7
# grades for the one student row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: clone an existing "Assignment / Grade / Comments / <sep>"
#        4-column block's formatting across the three new groups (AC:AN).
#        Q1:T16 is exactly one such block (separator + assignment/grade/comments
#        columns), and 12 columns = 3 x 4, so the paste tiles automatically.
$ws.Range("Q1:T16").Copy() | Out-Null
$ws.Range("AC1:AN16").PasteSpecial(-4122) | Out-Null

# --- 2. New header row labels for the 3 freshly-formatted groups.
$ws.Range("AD1").Value = "Assignment"
$ws.Range("AE1").Value = "Grade"
$ws.Range("AF1").Value = "Comments"
$ws.Range("AH1").Value = "Assignment"
$ws.Range("AI1").Value = "Grade"
$ws.Range("AJ1").Value = "Comments"
$ws.Range("AL1").Value = "Assignment"
$ws.Range("AM1").Value = "Grade"
$ws.Range("AN1").Value = "Comments"

# --- 3. Grade data for Sharon Vishnivetsky (row 2).
$ws.Range("B2").Value = 1

# HW5 (R:T) - was blank, now "4 UE" / 0 / "Not submitted"
# (write T2 before R2 so new shared-string entries land in the same order
# the original workbook's sharedStrings table uses)
$ws.Range("T2").Value = "Not submitted"
$ws.Range("R2").Value = "4 UE"
$ws.Range("S2").Value = 0

# HW6 (V:X) - was blank, now "5 UE" / 0 / "Not submitted"
$ws.Range("V2").Value = "5 UE"
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = "Not submitted"

# HW7 (Z:AB) - was blank, now "6 UE" / 0 / "Not submitted"
$ws.Range("Z2").Value = "6 UE"
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = "Not submitted"

# HW8 (AD:AF) - brand-new group, "7 CPP" / 63 / "my comments in the code"
$ws.Range("AD2").Value = "7 CPP"
$ws.Range("AE2").Value = 63
$ws.Range("AF2").Value = "my comments in the code"

# --- 4. Cosmetic: header row wraps to two lines, data row grew (long
#        comment), and the active selection moved while scrolled right.
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 60
$ws.Range("Z8").Select() | Out-Null
